# Housing Price Prediction Kaggle Competition
# Update the "Survived" (column B) predictions for a set of passenger rows
# from 0 to 1 on the gender_submission sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(6,20,26,28,50,54,55,61,68,71,76,94,98,102,114,116,124,143,152,155,158,161,170,181,186,210,216,220,222,226,240,241,242,265,274,316,326,345,349,352,366,370,373,377,399,402,404,411,413,416)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 2).Value = 1
}
